$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-08 Monday", "2024-04-09 Tuesday"),
    @("390×3=", "149×7="),
    @("780×9=", "731×4="),
    @("789×6=", "670×9="),
    @("844×3=", "659×9="),
    @("398×7=", "982×5="),
    @("927×4=", "760×3="),
    @("947×3=", "840×8="),
    @("821×4=", "162×3="),
    @("302×4=", "556×6="),
    @("128×4=", "300×6="),
    @("775×8=", "562×9="),
    @("381×2=", "603×7="),
    @("896×7=", "829×3="),
    @("268×8=", "148×8="),
    @("116×8=", "783×8="),
    @("499×9=", "760×8="),
    @("132×2=", "477×7="),
    @("762×6=", "455×6="),
    @("764×3=", "192×7="),
    @("225×5=", "208×6="),
    @("221×6=", "471×6="),
    @("574×9=", "574×8="),
    @("697×2=", "508×7="),
    @("930×9=", "880×6="),
    @("760×5=", "108×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
